$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'enden',
    'messen',
    'wenden',
    'fischen',
    'zünden',
    'seufzen',
    'dürfen',
    'rasen',
    'wohnen',
    'formen',
    'zielen',
    'bauen',
    'töten',
    'hauen',
    'brauen',
    'plaudern',
    'sprengen',
    'fallen',
    'segnen',
    'knarren',
    'mauern',
    'kehren',
    'hindern',
    'nutzen',
    'fühlen',
    'filmen',
    'mögen',
    'locken',
    'flehen',
    'geben',
    'kratzen',
    'bluten',
    'schwingen',
    'siegen',
    'fällen',
    'doppeln',
    'scheinen',
    'schultern',
    'flüchten',
    'buchen',
    'wachsen',
    'fahren',
    'grüßen',
    'brauchen',
    'sagen',
    'bitten',
    'quälen',
    'spinnen',
    'sperren',
    'nerven',
    'stehlen',
    'schauen',
    'warnen',
    'runden',
    'hupen',
    'loben',
    'rühren',
    'heben',
    'schlucken',
    'graben',
    'dringen',
    'opfern',
    'altern',
    'parken',
    'münzen',
    'folgen',
    'helfen',
    'jagen',
    'schulden',
    'atmen',
    'lockern',
    'sterben',
    'beten',
    'klappen',
    'bellen',
    'ziehen',
    'sichern',
    'rechnen',
    'suchen',
    'testen',
    'backen',
    'gelten',
    'trennen',
    'herzen',
    'träumen',
    'scheitern',
    'bergen',
    'ärgern',
    'freuen',
    'rufen',
    'trotzen',
    'wecken',
    'tollen',
    'schrecken',
    'zeigen',
    'heulen',
    'hören',
    'lügen',
    'warten',
    'werfen',
    'trauen',
    'wählen',
    'spielen',
    'merken',
    'binden',
    'schwören',
    'boxen',
    'fließen',
    'wärmen',
    'decken',
    'meinen',
    'schmecken',
    'schlagen',
    'tropfen',
    'zögern',
    'wechseln',
    'räumen',
    'führen',
    'kaufen',
    'ehren',
    'leihen',
    'dienen',
    'biegen',
    'greifen',
    'liefern',
    'äußern',
    'heilen',
    'wollen',
    'klingen',
    'schließen',
    'saufen',
    'wundern',
    'sammeln',
    'feiern',
    'gründen',
    'machen',
    'wehtun',
    'jubeln',
    'kichern',
    'lohnen',
    'wüten',
    'schalten',
    'pfeifen',
    'sorgen',
    'klettern',
    'lassen',
    'stammen',
    'achten',
    'schenken',
    'pressen',
    'erben',
    'krachen',
    'reizen',
    'fassen',
    'zeugen',
    'treiben',
    'weinen',
    'arten',
    'fügen',
    'betteln',
    'flüstern',
    'kümmern',
    'schützen',
    'lesen',
    'schwächen',
    'werden',
    'starten',
    'ändern',
    'drehen',
    'kosten',
    'pflanzen',
    'löschen',
    'weigern',
    'irren',
    'schreiten',
    'stecken',
    'kriegen',
    'schneiden',
    'knien',
    'malen',
    'liegen',
    'fesseln',
    'fangen',
    'planen',
    'sinken',
    'wirken',
    'weichen',
    'spüren',
    'prügeln',
    'stimmen',
    'streichen',
    'platzen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

